$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Bug fix: the scraper was pulling the wrong columns, so the 2006 and 2007
# tournament rows were missing and the 2019 row was never added. Insert the
# missing 2006/2007 rows near the top of the table and append the 2019 row.
# ---------------------------------------------------------------------------

# Insert 6 new rows right after the header for the 2006 and 2007 seasons,
# pushing the existing 2008-2018 data (and the trailing blank rows) down.
$ws.Rows("2:7").Insert()

# --- Column A (Year) -------------------------------------------------------
# Enter the numeric year values first, then apply the text format used by the
# rest of the table (matches the existing quirk where Year is typed as a
# number but displayed with the worksheet's text format).
$ws.Range("A2").Value2 = 2006
$ws.Range("A3").Value2 = 2006
$ws.Range("A4").Value2 = 2006
$ws.Range("A5").Value2 = 2007
$ws.Range("A6").Value2 = 2007
$ws.Range("A7").Value2 = 2007
$ws.Range("A2:A7").NumberFormat = "@"

# --- Column B (Round) -------------------------------------------------------
$ws.Range("B2:B7").Style = "Normal"
$ws.Range("B2").Value2 = "quarter-finals"
$ws.Range("B3").Value2 = "semi-finals"
$ws.Range("B4").Value2 = "finals"
$ws.Range("B5").Value2 = "quarter-finals"
$ws.Range("B6").Value2 = "semi-finals"
$ws.Range("B7").Value2 = "finals"

# --- Columns C & D (Start / End) -------------------------------------------
$ws.Range("C2:D7").NumberFormat = "@"
$ws.Range("C2").Value2 = "2006-04-21"
$ws.Range("D2").Value2 = "2006-05-04"
$ws.Range("C3").Value2 = "2006-04-21"
$ws.Range("D3").Value2 = "2006-05-18"
$ws.Range("C4").Value2 = "2006-04-21"
$ws.Range("D4").Value2 = "2006-06-04"
$ws.Range("C5").Value2 = "2007-04-11"
$ws.Range("D5").Value2 = "2007-04-24"
$ws.Range("C6").Value2 = "2007-04-11"
$ws.Range("D6").Value2 = "2007-05-09"
$ws.Range("C7").Value2 = "2007-04-11"
$ws.Range("D7").Value2 = "2006-05-27"

# ---------------------------------------------------------------------------
# Add the missing 2019 row. After the insert above, row 41 is the first of
# the pre-existing blank rows at the bottom of the sheet, already carrying
# the table's text style - just fill in its values.
# ---------------------------------------------------------------------------
$ws.Range("A41").Value2 = "2019"
$ws.Range("B41").Value2 = "quarter-finals"
$ws.Range("C41").Value2 = "2019-04-10"
$ws.Range("D41").Value2 = "2019-04-24"

# Update the view: scroll down a bit and select D42, the cell right after the
# newly-added data.
$ws.Activate()
$ws.Range("D42").Select()
